# Updated the main page
#
# Inserts a new "REPO" slide (Title + Content layout) right before the
# closing "Feedback" slide, linking to the workshop's GitHub repository.

$p = $ppt.ActivePresentation

# The deck currently ends with: ... , Agenda (afternoon), Feedback
# We insert the new slide at position 5, pushing "Feedback" to position 6.
$feedbackIndex = $p.Slides.Count
$newSlide = $p.Slides.Add($feedbackIndex, 2)

# Title placeholder -> "REPO"
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "REPO"

# Content placeholder -> hyperlinked GitHub URL
$repoUrl = "https://github.com/mufajjul/aml-govsec2020-workshop"
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = $repoUrl

$bodyRun = $body.Characters(1, $body.Length)
$bodyRun.ActionSettings.Item(1).Hyperlink.Address = $repoUrl
$bodyRun.LanguageID = "en-GB"
